$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 25.1800000000005
$ws.Range("H2").Value = 0.6236307111900747
$ws.Range("I2").Value = 0.6236307111900747
$ws.Range("L2").Value = 3.536826012802427
$ws.Range("M2").Value = "[-6.035819834482176, 13.10947186008703]"
$ws.Range("N2").Value = 0.4606510147127172
$ws.Range("O2").Value = 0.4606510147127172
$ws.Range("P2").Value = -2.666737307551696
$ws.Range("Q2").Value = "[-5.805185223750509, 0.47171060864711567]"
$ws.Range("R2").Value = 0.09389603997900875
$ws.Range("S2").Value = 0.09389603997900875
$ws.Range("T2").Value = 12.62425206280053
$ws.Range("U2").Value = "[7.581233157310514, 17.66727096829054]"
$ws.Range("V2").Value = [double]"7.999352174214991e-06"
$ws.Range("W2").Value = [double]"7.999352174214991e-06"
$ws.Range("X2").Value = 10.68700700700722
$ws.Range("Y2").Value = -1.890390390390426
$ws.Range("Z2").Value = 23.26440440440486

# Row 3
$ws.Range("F3").Value = 25.1800000000005
$ws.Range("H3").Value = 0.4550530020617787
$ws.Range("I3").Value = 0.4550530020617787
$ws.Range("L3").Value = 4.082413603826061
$ws.Range("M3").Value = "[-4.4602142479561815, 12.625041455608304]"
$ws.Range("N3").Value = 0.3409351045504181
$ws.Range("O3").Value = 0.3409351045504181
$ws.Range("P3").Value = -2.490632013656773
$ws.Range("Q3").Value = "[-5.622790455073624, 0.6415264277600787]"
$ws.Range("R3").Value = 0.1162469533482031
$ws.Range("S3").Value = 0.1162469533482031
$ws.Range("T3").Value = 10.61097210360921
$ws.Range("U3").Value = "[6.093544122899601, 15.128400084318827]"
$ws.Range("V3").Value = [double]"2.240615799431467e-05"
$ws.Range("W3").Value = [double]"2.240615799431467e-05"
$ws.Range("X3").Value = 9.981261261261457
$ws.Range("Y3").Value = -2.570930930930986
$ws.Range("Z3").Value = 22.5334534534539

# Row 4
$ws.Range("F4").Value = 25.1800000000005
$ws.Range("H4").Value = 0.2219720947756371
$ws.Range("I4").Value = 0.2219720947756371
$ws.Range("L4").Value = 5.311032761899888
$ws.Range("M4").Value = "[-2.4423849988912743, 13.06445052269105]"
$ws.Range("N4").Value = 0.1745131737459444
$ws.Range("O4").Value = 0.1745131737459444
$ws.Range("P4").Value = -1.74847398938531
$ws.Range("Q4").Value = "[-4.7925797838547, 1.2956318050840796]"
$ws.Range("R4").Value = 0.2534335077141101
$ws.Range("S4").Value = 0.2534335077141101
$ws.Range("T4").Value = 14.42189596725099
$ws.Range("U4").Value = "[10.057712494392115, 18.786079440109866]"
$ws.Range("V4").Value = [double]"3.306285578652535e-08"
$ws.Range("W4").Value = [double]"3.306285578652535e-08"
$ws.Range("X4").Value = 7.007047047047184
$ws.Range("Y4").Value = -5.19227227227238
$ws.Range("Z4").Value = 19.20636636636675

# Row 5
$ws.Range("F5").Value = 25.1800000000005
$ws.Range("H5").Value = 0.6308737645698584
$ws.Range("I5").Value = 0.6308737645698584
$ws.Range("L5").Value = 3.063423554404779
$ws.Range("M5").Value = "[-5.147819971164046, 11.274667079973604]"
$ws.Range("N5").Value = 0.4563149008163125
$ws.Range("O5").Value = 0.4563149008163125
$ws.Range("P5").Value = -3.044105794469389
$ws.Range("Q5").Value = "[-6.138527387194471, 0.0503157982556921]"
$ws.Range("R5").Value = 0.05368142950139143
$ws.Range("S5").Value = 0.05368142950139143
$ws.Range("T5").Value = 10.97579262488651
$ws.Range("U5").Value = "[6.566727151971358, 15.384858097801652]"
$ws.Range("V5").Value = [double]"8.78445177354692e-06"
$ws.Range("W5").Value = [double]"8.78445177354692e-06"
$ws.Range("X5").Value = 12.19931931931956
$ws.Range("Y5").Value = -0.2016416416416433
$ws.Range("Z5").Value = 24.60028028028077

# Row 6
$ws.Range("F6").Value = 25.1800000000005
$ws.Range("H6").Value = 0.09465078282533401
$ws.Range("I6").Value = 0.09465078282533401
$ws.Range("L6").Value = 9.887224196821636
$ws.Range("M6").Value = "[-1.9578969437203355, 21.732345337363608]"
$ws.Range("N6").Value = 0.09965511262188187
$ws.Range("O6").Value = 0.09965511262188187
$ws.Range("P6").Value = -2.553526761476389
$ws.Range("Q6").Value = "[-5.666816778547354, 0.5597632555945768]"
$ws.Range("R6").Value = 0.1055023576923992
$ws.Range("S6").Value = 0.1055023576923992
$ws.Range("T6").Value = 15.56947677681293
$ws.Range("U6").Value = "[9.365469108791762, 21.773484444834104]"
$ws.Range("V6").Value = [double]"7.669217862771305e-06"
$ws.Range("W6").Value = [double]"7.669217862771305e-06"
$ws.Range("X6").Value = 10.23331331331352
$ws.Range("Y6").Value = -2.243263263263305
$ws.Range("Z6").Value = 22.70988988989034

# Row 7
$ws.Range("F7").Value = 25.1800000000005
$ws.Range("H7").Value = 0.2403782294450264
$ws.Range("I7").Value = 0.2403782294450264
$ws.Range("L7").Value = 6.596920177418399
$ws.Range("M7").Value = "[-3.2788807669549804, 16.47272112179178]"
$ws.Range("N7").Value = 0.1852375100102228
$ws.Range("O7").Value = 0.1852375100102228
$ws.Range("P7").Value = 2.798816277972889
$ws.Range("Q7").Value = "[-0.1823947686768852, 5.7800273246226626]"
$ws.Range("R7").Value = 0.0650893667522221
$ws.Range("S7").Value = 0.0650893667522221
$ws.Range("T7").Value = 13.24637390246868
$ws.Range("U7").Value = "[7.83231813410714, 18.66042967083021]"
$ws.Range("V7").Value = [double]"1.169326593708497e-05"
$ws.Range("W7").Value = [double]"1.169326593708497e-05"
$ws.Range("X7").Value = 13.96368368368396
$ws.Range("Y7").Value = 2.016416416416458
$ws.Range("Z7").Value = 25.91095095095146

# Row 8
$ws.Range("F8").Value = 25.1800000000005
$ws.Range("H8").Value = 0.2464275404997598
$ws.Range("I8").Value = 0.2464275404997598
$ws.Range("L8").Value = 5.650830639623713
$ws.Range("M8").Value = "[-3.1588817919674987, 14.460543071214925]"
$ws.Range("N8").Value = 0.2029840317219684
$ws.Range("O8").Value = 0.2029840317219684
$ws.Range("P8").Value = 2.698184681461504
$ws.Range("Q8").Value = "[-0.3522105877898465, 5.7485799507128545]"
$ws.Range("R8").Value = 0.08157261435804597
$ws.Range("S8").Value = 0.08157261435804597
$ws.Range("T8").Value = 12.92610951270433
$ws.Range("U8").Value = "[8.221476318957695, 17.630742706450963]"
$ws.Range("V8").Value = [double]"1.527936318579393e-06"
$ws.Range("W8").Value = [double]"1.527936318579393e-06"
$ws.Range("X8").Value = 14.36696696696725
$ws.Range("Y8").Value = 2.142442442442483
$ws.Range("Z8").Value = 26.59149149149201

# Row 9
$ws.Range("F9").Value = 25.1800000000005
$ws.Range("H9").Value = 0.8864652765304512
$ws.Range("I9").Value = 0.8864652765304512
$ws.Range("L9").Value = 1.727849151956206
$ws.Range("M9").Value = "[-8.261735230392816, 11.71743353430523]"
$ws.Range("N9").Value = 0.7291872549435507
$ws.Range("O9").Value = 0.7291872549435507
$ws.Range("P9").Value = 1.83023716155081
$ws.Range("Q9").Value = "[-1.3082107546480017, 4.9686850777496225]"
$ws.Range("R9").Value = 0.2463501070892007
$ws.Range("S9").Value = 0.2463501070892007
$ws.Range("T9").Value = 15.66192074982151
$ws.Range("U9").Value = "[10.531940277943733, 20.79190122169929]"
$ws.Range("V9").Value = [double]"1.875058051048484e-07"
$ws.Range("W9").Value = [double]"1.875058051048484e-07"
$ws.Range("X9").Value = 17.84528528528564
$ws.Range("Y9").Value = 5.267887887887992
$ws.Range("Z9").Value = 30.42268268268328

# Row 10
$ws.Range("F10").Value = 24.32000000000036
$ws.Range("H10").Value = 0.7598435546755605
$ws.Range("I10").Value = 0.7598435546755605
$ws.Range("L10").Value = 2.515459190501703
$ws.Range("M10").Value = "[-6.142888628157877, 11.173807009161283]"
$ws.Range("N10").Value = 0.5613726544336739
$ws.Range("O10").Value = 0.5613726544336739
$ws.Range("P10").Value = 1.025184389459733
$ws.Range("Q10").Value = "[-2.11326352673908, 4.163632305658546]"
$ws.Range("R10").Value = 0.5139454326158006
$ws.Range("S10").Value = 0.5139454326158006
$ws.Range("T10").Value = 14.1661754386351
$ws.Range("U10").Value = "[9.340264246550621, 18.99208663071957]"
$ws.Range("V10").Value = [double]"4.213406612496584e-07"
$ws.Range("W10").Value = [double]"4.213406612496584e-07"
$ws.Range("X10").Value = 20.35187187187217
$ws.Range("Y10").Value = 8.204044044044158
$ws.Range("Z10").Value = 32.49969969970018

# Row 11
$ws.Range("F11").Value = 24.32000000000036
$ws.Range("H11").Value = 0.7127835887910903
$ws.Range("I11").Value = 0.7127835887910903
$ws.Range("L11").Value = 2.859631026898689
$ws.Range("M11").Value = "[-6.456739250421951, 12.17600130421933]"
$ws.Range("N11").Value = 0.5395451120613952
$ws.Range("O11").Value = 0.5395451120613952
$ws.Range("P11").Value = 2.371131992799504
$ws.Range("Q11").Value = "[-0.767315923399309, 5.509579908998316]"
$ws.Range("R11").Value = 0.1350865244297508
$ws.Range("S11").Value = 0.1350865244297508
$ws.Range("T11").Value = 12.66323491981698
$ws.Range("U11").Value = "[7.725961855921042, 17.600507983712916]"
$ws.Range("V11").Value = [double]"5.285898810214462e-06"
$ws.Range("W11").Value = [double]"5.285898810214462e-06"
$ws.Range("X11").Value = 15.14218218218241
$ws.Range("Y11").Value = 2.994354354354398
$ws.Range("Z11").Value = 27.29001001001042

# Row 12
$ws.Range("F12").Value = 24.32000000000036
$ws.Range("H12").Value = 0.9171674963803411
$ws.Range("I12").Value = 0.9171674963803411
$ws.Range("L12").Value = 1.073299089347797
$ws.Range("M12").Value = "[-5.793065366035194, 7.939663544730788]"
$ws.Range("N12").Value = 0.7543461944637582
$ws.Range("O12").Value = 0.7543461944637582
$ws.Range("P12").Value = 1.754763464167272
$ws.Range("Q12").Value = "[-1.377394977249578, 4.886921905584122]"
$ws.Range("R12").Value = 0.2651381843710152
$ws.Range("S12").Value = 0.2651381843710152
$ws.Range("T12").Value = 12.37384290634153
$ws.Range("U12").Value = "[8.678652868038313, 16.069032944644754]"
$ws.Range("V12").Value = [double]"2.440640134437899e-08"
$ws.Range("W12").Value = [double]"2.440640134437899e-08"
$ws.Range("X12").Value = 17.52792792792819
$ws.Range("Y12").Value = 5.404444444444529
$ws.Range("Z12").Value = 29.65141141141185

# Row 13
$ws.Range("F13").Value = 24.32000000000036
$ws.Range("H13").Value = 0.5927960663844726
$ws.Range("I13").Value = 0.5927960663844726
$ws.Range("L13").Value = 3.089615601787088
$ws.Range("M13").Value = "[-5.068274646061205, 11.247505849635381]"
$ws.Range("N13").Value = 0.4495654387705987
$ws.Range("O13").Value = 0.4495654387705987
$ws.Range("P13").Value = 2.157289850212811
$ws.Range("Q13").Value = "[-0.9622896416401172, 5.276869342065739]"
$ws.Range("R13").Value = 0.1705185015541737
$ws.Range("S13").Value = 0.1705185015541737
$ws.Range("T13").Value = 10.68921372749509
$ws.Range("U13").Value = "[6.394646256632816, 14.983781198357372]"
$ws.Range("V13").Value = [double]"8.805767878650173e-06"
$ws.Range("W13").Value = [double]"8.805767878650173e-06"
$ws.Range("X13").Value = 15.96988988989013
$ws.Range("Y13").Value = 3.895095095095153
$ws.Range("Z13").Value = 28.0446846846851

# Row 14
$ws.Range("F14").Value = 24.32000000000036
$ws.Range("H14").Value = 0.4021563475787338
$ws.Range("I14").Value = 0.4021563475787338
$ws.Range("L14").Value = 3.66790043972779
$ws.Range("M14").Value = "[-3.7289795097088465, 11.064780389164426]"
$ws.Range("N14").Value = 0.3232629766519211
$ws.Range("O14").Value = 0.3232629766519211
$ws.Range("P14").Value = 1.188710733790733
$ws.Range("Q14").Value = "[-1.9497371824080796, 4.327158649989545]"
$ws.Range("R14").Value = 0.4495294845627864
$ws.Range("S14").Value = 0.4495294845627864
$ws.Range("T14").Value = 11.10639268496179
$ws.Range("U14").Value = "[7.260674168317919, 14.952111201605664]"
$ws.Range("V14").Value = [double]"5.837979442091523e-07"
$ws.Range("W14").Value = [double]"5.837979442091523e-07"
$ws.Range("X14").Value = 19.71891891891921
$ws.Range("Y14").Value = 7.571091091091208
$ws.Range("Z14").Value = 31.86674674674722

# Row 15
$ws.Range("B15").Value = 1
$ws.Range("F15").Value = 24.32000000000036
$ws.Range("H15").Value = 0.01939542259931271
$ws.Range("I15").Value = 0.01939542259931271
$ws.Range("L15").Value = 9.162391083448648
$ws.Range("M15").Value = "[2.039202783332417, 16.28557938356488]"
$ws.Range("N15").Value = 0.01286447415408731
$ws.Range("O15").Value = 0.01286447415408731
$ws.Range("P15").Value = 1.566079220708426
$ws.Range("Q15").Value = "[0.4465527095192705, 2.6856057318975814]"
$ws.Range("R15").Value = 0.007164936334548511
$ws.Range("S15").Value = 0.007164936334548511
$ws.Range("T15").Value = 12.05521947268288
$ws.Range("U15").Value = "[7.581315375979987, 16.52912356938578]"
$ws.Range("V15").Value = [double]"2.192650618848546e-06"
$ws.Range("W15").Value = [double]"2.192650618848546e-06"
$ws.Range("X15").Value = 18.25825825825853
$ws.Range("Y15").Value = 13.92496496496517
$ws.Range("Z15").Value = 22.59155155155189
